$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) for rows 2-11
$ws.Range("C2").Value = 0.9310122474908051
$ws.Range("D2").Value = 0.3584111998068638

$ws.Range("C3").Value = 2.387288937183493
$ws.Range("D3").Value = 0.0226774307765456

$ws.Range("C4").Value = 1.01172709995921
$ws.Range("D4").Value = 0.3188147116121174

$ws.Range("C5").Value = 3.390954101638172
$ws.Range("D5").Value = 0.001779768824629269

$ws.Range("C6").Value = 1.248828773388428
$ws.Range("D6").Value = 0.2202581191770054

$ws.Range("C7").Value = 0.3179953770968922
$ws.Range("D7").Value = 0.7524334188045052

$ws.Range("C8").Value = 2.827373761711713
$ws.Range("D8").Value = 0.007809069868069729

$ws.Range("C9").Value = -0.9768400401535045
$ws.Range("D9").Value = 0.3355479016353935

$ws.Range("C10").Value = 1.082125622388216
$ws.Range("D10").Value = 0.2868120999653208

$ws.Range("C11").Value = 1.886639140527811
$ws.Range("D11").Value = 0.06777275666478855

# Row 11 "Mejor" column changes from "Sí" to "No"
$ws.Range("G11").Value = "No"
